{"js": "// Add a \"Meta description\" paragraph right after the document title, remove the\n// duplicate title paragraph near the end of the document, and replace the text of\n// the trailing italic paragraph with the new \"Prompt: ...\" image-generation prompt.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst titleText = \"Play Book of Darkness for Free - Betsoft Slot Review\";\nconst oldMetaText =\n  \"Try for free Book of Darkness, the high-volatility slot game where you can win big. Our review highlights its engaging gameplay, storyline, and special features.\";\nconst newPromptText =\n  \"Prompt: Create a fun and engaging feature image for Book of Darkness that captures the thrilling concept of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses, surrounded by fire and magic symbols. Make sure to include the iconic Book of Darkness in the image as well. Overall, the image should showcase the power struggle between the good and evil forces in the game, with the Maya warrior as the hero fighting against the evil magician. The image should be vibrant and eye-catching, with bold colors and dynamic imagery that draws in the player's attention.\";\n\n// --- Step 1: insert the new \"Meta description\" paragraph right after the title ---\nlet titleParaIndex = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === titleText) {\n    titleParaIndex = i;\n    break;\n  }\n}\nconst titlePara = paragraphs.items[titleParaIndex];\n\nconst metaParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  \"<w:r/>\" +\n  \"<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">: Try for free Book of Darkness, the high-volatility slot game where you can win big. Our review highlights its engaging gameplay, storyline, and special features.</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:sectPr/>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ntitlePara.getRange().insertOoxml(metaParagraphOoxml, Word.InsertLocation.after);\nawait context.sync();\n\n// --- Step 2: delete the duplicate bold title paragraph further down the document ---\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet duplicateTitleIndex = -1;\nfor (let i = 1; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === titleText) {\n    duplicateTitleIndex = i;\n    break;\n  }\n}\nif (duplicateTitleIndex !== -1) {\n  paragraphs.items[duplicateTitleIndex].delete();\n  await context.sync();\n}\n\n// --- Step 3: replace the text of the trailing italic paragraph with the new prompt ---\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet promptParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === oldMetaText) {\n    promptParaIndex = i;\n    break;\n  }\n}\nif (promptParaIndex !== -1) {\n  paragraphs.items[promptParaIndex]\n    .getRange()\n    .insertText(newPromptText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Add a \"Meta description\" paragraph right after the document title, remove the\n# duplicate title paragraph near the end of the document, and replace the text of\n# the trailing italic paragraph with the new \"Prompt: ...\" image-generation prompt.\n\n$d = $word.ActiveDocument\n\n$titleText = \"Play Book of Darkness for Free - Betsoft Slot Review\"\n$oldPromptText = \"Try for free Book of Darkness, the high-volatility slot game where you can win big. Our review highlights its engaging gameplay, storyline, and special features.\"\n$newPromptText = \"Prompt: Create a fun and engaging feature image for Book of Darkness that captures the thrilling concept of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses, surrounded by fire and magic symbols. Make sure to include the iconic Book of Darkness in the image as well. Overall, the image should showcase the power struggle between the good and evil forces in the game, with the Maya warrior as the hero fighting against the evil magician. The image should be vibrant and eye-catching, with bold colors and dynamic imagery that draws in the player's attention.\"\n\n# --- Step 1: insert the new \"Meta description\" paragraph right after the title ---\n# NOTE: this engine's Range.InsertXML(), when given a whole <w:p>...</w:p> block and\n# a *collapsed* (zero-length) range sitting exactly on a paragraph boundary, replaces\n# the paragraph adjacent to that boundary instead of inserting a sibling next to it.\n# To safely insert a new paragraph right after the title without clobbering either\n# neighbour, target the title paragraph's own (non-collapsed) range and feed\n# InsertXML an OOXML fragment that reproduces the title paragraph followed by the\n# new paragraph - i.e. a same-content \"replace\" that also appends the new paragraph.\n$titleParaIndex = 1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -eq ($titleText + \"`r\")) {\n        $titleParaIndex = $i\n        break\n    }\n}\n$titlePara = $d.Paragraphs.Item($titleParaIndex)\n$titleRange = $titlePara.Range\n$titleParaOwnText = $titlePara.Range.Text.TrimEnd([char]13)\n# Defensive XML-escaping in case the title ever contains markup-sensitive characters.\n$titleParaOwnTextXml = $titleParaOwnText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n\n$metaParagraphOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:t>' + $titleParaOwnTextXml + '</w:t></w:r></w:p>' +\n  '<w:p>' +\n  '<w:r/>' +\n  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">: Try for free Book of Darkness, the high-volatility slot game where you can win big. Our review highlights its engaging gameplay, storyline, and special features.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:sectPr/>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n$titleRange.InsertXML($metaParagraphOoxml)\n\n# --- Step 2: delete the duplicate bold title paragraph further down the document ---\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 2; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq ($titleText + \"`r\")) {\n        $p.Range.Delete()\n    }\n}\n\n# --- Step 3: replace the text of the trailing italic paragraph with the new prompt ---\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq ($oldPromptText + \"`r\")) {\n        $r = $d.Range($p.Range.Start, $p.Range.End - 1)\n        $r.Text = $newPromptText\n    }\n}\n"}
